# Update Avanzamento sheet rows 66-131 (Ore lavorate / Produzione refresh)
# and restore the sheet view (scroll position + selection) as recorded in the
# workbook when it was last saved from Streamlit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> @{ C = <Ore lavorate>; D = <Produzione> }
# (only the rows whose Ore lavorate / Produzione figures were refreshed)
$rowData = @{
  66 = @{ D=29.237926973373401 }
  67 = @{ D=41.853499128812601 }
  68 = @{ C=168; D=26.775168090786401 }
  69 = @{ D=48.311080659107397 }
  70 = @{ C=176; D=28.5732775665399 }
  71 = @{ C=120; D=40.093497186019398 }
  72 = @{ D=31.780536585365901 }
  73 = @{ C=128; D=33.741710735250201 }
  74 = @{ D=59.2818476108068 }
  75 = @{ C=184; D=60.856685425766599 }
  76 = @{ D=52.144468873261999 }
  77 = @{ D=44.972536565901301 }
  78 = @{ C=168; D=28.5732775665399 }
  79 = @{ D=35.665015290806799 }
  81 = @{ D=43.636136627907 }
  82 = @{ D=25.857861904483901 }
  83 = @{ D=31.253423015873 }
  84 = @{ C=175; D=44.031594599410397 }
  85 = @{ D=30.190300707236801 }
  86 = @{ C=144; D=49.901317293941098 }
  87 = @{ C=144; D=40.917800928913003 }
  88 = @{ C=176; D=35.802042386979601 }
  89 = @{ C=132; D=49.698799846698897 }
  90 = @{ C=160; D=92.0665249242561 }
  91 = @{ C=104; D=32.306667307692301 }
  92 = @{ C=158; D=25.417425653784399 }
  93 = @{ C=140; D=47.438137801035097 }
  94 = @{ C=176; D=44.165378198286298 }
  96 = @{ C=168; D=43.721122756013301 }
  97 = @{ C=168; D=81.029703389830502 }
  99 = @{ C=172; D=50.860851520567998 }
  100 = @{ C=168; D=41.517803674055799 }
  101 = @{ C=144; D=93.5786039218506 }
  102 = @{ D=25.9888573411728 }
  103 = @{ C=156; D=27.627710261205099 }
  104 = @{ D=38.4600450131976 }
  106 = @{ C=112; D=73.448572987089804 }
  107 = @{ D=26.673684976789598 }
  108 = @{ C=160; D=35.714481249999999 }
  109 = @{ C=176; D=81.029703389830502 }
  110 = @{ C=104; D=81.029703389830502 }
  111 = @{ C=206; D=42.041652962155702 }
  112 = @{ C=166; D=35.0972051363348 }
  113 = @{ C=182; D=41.72819933281 }
  114 = @{ C=198; D=33.8435735713743 }
  115 = @{ D=32.458925520833397 }
  116 = @{ C=157; D=39.0950206798762 }
  117 = @{ D=32.149164236111098 }
  118 = @{ C=176; D=31.474488636363599 }
  119 = @{ C=104; D=29.6587761470985 }
  120 = @{ C=176; D=28.5732775665399 }
  121 = @{ C=112; D=68.226912149263498 }
  122 = @{ C=176; D=40.862349444089901 }
  123 = @{ C=168; D=41.589169642857101 }
  124 = @{ D=30.136456849049701 }
  125 = @{ D=47.735689272827202 }
  126 = @{ D=70.325269855422803 }
  128 = @{ C=176; D=105.384265957447 }
  129 = @{ C=96; D=27.6790264900662 }
  130 = @{ C=112; D=35.714481249999999 }
  131 = @{ C=104; D=33.073374999999999 }
}

foreach ($row in $rowData.Keys) {
    $entry = $rowData[$row]
    if ($entry.ContainsKey("C")) {
        $ws.Cells.Item([int]$row, 3).Value = $entry["C"]
    }
    if ($entry.ContainsKey("D")) {
        $ws.Cells.Item([int]$row, 4).Value = $entry["D"]
    }
}

# Restore the view: scrolled so row 48 is at the top, with G66 selected
# (matches the sheetView/selection recorded in the saved workbook).
$win = $excel.ActiveWindow
$win.ScrollRow = 48
$win.ScrollColumn = 1
$ws.Range("G66").Select()
